$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.661.86"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.291.41"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.92"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.56"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.290.31"
$ws.Range("E9").Value = "  -5.33%  "
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.53"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.77"
$ws.Range("E14").Value = "  -4.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.698.92"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.637.56"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.261.19"
$ws.Range("E18").Value = "  -6.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.65"
$ws.Range("E19").Value = "  -5.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.30"
$ws.Range("E21").Value = "  -3.61%  "
$ws.Range("E22").Value = "  -4.13%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.02"
$ws.Range("E24").Value = "  -3.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.08"
$ws.Range("E27").Value = "  -7.04%  "
$ws.Range("E28").Value = "  -6.05%  "
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.56"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0726"
$ws.Range("E31").Value = "  -6.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.79"
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.383"
$ws.Range("E34").Value = "  -5.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.79"
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("E38").Value = "  -6.17%  "
$ws.Range("E39").Value = "  -6.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.95"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("E41").Value = "  -5.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "299.33"
$ws.Range("E42").Value = "  -8.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.20"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("E44").Value = "  -5.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0953"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0500"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.556"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.52"
$ws.Range("E48").Value = "  -7.42%  "
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.64"
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("E51").Value = "  -0.43%  "
